$d = $word.ActiveDocument

# Paragraph 1 ("All career schools") is kept as-is.
# Every paragraph after it (the five external-link hyperlinks, the blank
# Hyperlink-styled paragraph, the "alberta" paragraph and the final
# academyoflearning hyperlink paragraph that carries the _GoBack bookmark)
# is removed, but the _GoBack bookmark itself must end up anchored right
# after the surviving run, inside paragraph 1.

$firstPara = $d.Paragraphs.Item(1)

# Delete everything from the end of paragraph 1 (i.e. right after its
# paragraph mark) through to the end of the document content. This also
# removes the original _GoBack bookmark along with the paragraphs that
# held it.
$deleteStart = $firstPara.Range.End
$deleteEnd = $d.Content.End
$d.Range($deleteStart, $deleteEnd).Delete()

# Temporarily insert a placeholder character right before the paragraph
# mark of (now the only) paragraph 1, so that the insertion point used for
# the bookmark is a genuine mid-paragraph position rather than the exact
# paragraph-end position (which gets normalized/widened by Bookmarks.Add
# to span the whole paragraph). This lets us get a truly collapsed
# bookmark positioned immediately after "All career schools".
$p1 = $d.Paragraphs.Item(1)
$endPos = $p1.Range.End - 1
$d.Range($endPos, $endPos).InsertAfter("X")

# Re-create the _GoBack bookmark collapsed right before the placeholder.
$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$bm = $d.Bookmarks.Item("_GoBack")

# Remove the placeholder character again (it now sits right after the
# bookmark), leaving the bookmark correctly collapsed at the end of the
# run text, before the paragraph mark.
$d.Range($bm.End, $bm.End + 1).Delete()
